$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.456.22"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").Value = "1.891.68"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.84"
$ws.Range("E5").Value = "  +0.89%  "
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4892"
$ws.Range("E7").Value = "  +0.37%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2930"
$ws.Range("E8").Value = "  +1.26%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06677"
$ws.Range("E9").Value = "  +0.26%  "
$ws.Range("D10").Value = "1.884.62"
$ws.Range("E10").Value = "  -0.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "17.08"
$ws.Range("E11").Value = "  +2.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07344"
$ws.Range("E12").Value = "  +1.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.140"
$ws.Range("E13").Value = "  +2.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.15"
$ws.Range("E14").Value = "  -1.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6640"
$ws.Range("E15").Value = "  +0.14%  "
$ws.Range("D16").Value = "30.422.63"
$ws.Range("E16").Value = "  -0.67%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.43"
$ws.Range("E17").Value = "  +3.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007817"
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("D20").Value = "2.120.02"
$ws.Range("E20").Value = "  -0.42%  "
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.290"
$ws.Range("E22").Value = "  +11.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "188.33"
$ws.Range("E23").Value = "  -1.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.152"
$ws.Range("E24").Value = "  +1.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.476"
$ws.Range("E25").Value = "  +1.97%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.74"
$ws.Range("E26").Value = "  +3.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.31"
$ws.Range("E27").Value = "  +0.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.929"
$ws.Range("E28").Value = "  +5.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.464"
$ws.Range("E29").Value = "  +4.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.347"
$ws.Range("E30").Value = "  +2.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09167"
$ws.Range("E31").Value = "  +1.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.089"
$ws.Range("E32").Value = "  +4.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05200"
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("E34").Value = "  +1.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.097"
$ws.Range("E35").Value = "  +1.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.719"
$ws.Range("E36").Value = "  +1.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01820"
$ws.Range("E37").Value = "  -0.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.674"
$ws.Range("E38").Value = "  +0.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9164"
$ws.Range("E39").Value = "  -0.67%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.033"
$ws.Range("E40").Value = "  -0.82%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4393"
$ws.Range("E41").Value = "  +0.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.942"
$ws.Range("E42").Value = "  +4.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "105.60"
$ws.Range("E43").Value = "  +1.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9929"
$ws.Range("E44").Value = "  -0.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1385"
$ws.Range("E45").Value = "  +3.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "68.14"
$ws.Range("E46").Value = "  +18.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.563"
$ws.Range("E47").Value = "  +3.89%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.952"
$ws.Range("E48").Value = "  +3.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.94"
$ws.Range("E49").Value = "  +5.25%  "
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3936"
$ws.Range("E51").Value = "  -3.89%  "
